$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = "Форма государственной статистической отчетности №10 «Отчет о психических и наркологических расстройствах»"
